$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete column E (the duplicate "Ki hieu" column that had no body data),
# shifting column F ("days") left into column E.
$ws.Range("E:E").Delete()

# Update the selection/view to match the post-edit state.
$ws.Range("E1:E1048576").Select()
